# Update the "Date:" / "Time:" stamp embedded in the OLS Regression Results
# text block that lives in cell B2 of every worksheet (one worksheet per
# backward-elimination step). The workbook was re-generated/re-saved the
# next day, so every occurrence of the old timestamp is replaced with the
# new one, leaving the rest of the (very long) text untouched.

$wb = $excel.ActiveWorkbook

$oldDate = "Sat, 28 Dec 2019"
$newDate = "Sun, 29 Dec 2019"
$oldTime = "21:00:00"
$newTime = "16:11:33"

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("B2")
    $text = $cell.Value()

    if ($text -eq $null) { continue }
    if (-not ($text.Contains($oldDate) -or $text.Contains($oldTime))) { continue }

    $row = $ws.Rows.Item(2)
    $origRowHeight = $row.RowHeight

    # Excel's COM layer hands back the cached text with normalized LF line
    # breaks even though the text is stored as CRLF in the file; restore the
    # CRLF endings before writing the value back so the round-tripped XML
    # stays byte-identical aside from the intended replacement.
    $text = $text.Replace("`r`n", "`n").Replace("`n", "`r`n")

    $text = $text.Replace($oldDate, $newDate)
    $text = $text.Replace($oldTime, $newTime)

    $cell.Value = $text

    # Re-assigning this wrapped, 60+ line cell makes the host recompute the
    # row's auto-fit height; restore the original (already-at-the-Excel-max
    # 409.5pt) height if the write nudged it away from its prior value so
    # the row stays exactly as it was.
    if ($row.RowHeight -ne $origRowHeight) {
        $row.RowHeight = $origRowHeight
    }
}
